$wb = $excel.ActiveWorkbook

# weibull
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.49972872840229
$ws.Range("C2").Value = 0.122057447093382
$ws.Range("B3").Value = 0.0391414385130606
$ws.Range("C3").Value = 0.0517908396172799

# lognormal
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 1.93563637872471
$ws.Range("C2").Value = 0.121436695299748
$ws.Range("B3").Value = -0.991804474991053
$ws.Range("C3").Value = 0.0583512054430059

# llogis
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.93060133863177
$ws.Range("C2").Value = 0.102811464864021
$ws.Range("B3").Value = 1.71540910960479
$ws.Range("C3").Value = 0.128615701649951

# gompertz
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.29339372378354
$ws.Range("C2").Value = 0.132301106530333
$ws.Range("B3").Value = -0.013784129933345
$ws.Range("C3").Value = 0.00673971626319614

# weibull cov
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0148980203909537
$ws.Range("B2").Value = -0.00257155016007909
$ws.Range("A3").Value = -0.00257155016007909
$ws.Range("B3").Value = 0.00268229106826281

# lognormal cov
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0147468709653237
$ws.Range("B2").Value = -0.00444446528904442
$ws.Range("A3").Value = -0.00444446528904442
$ws.Range("B3").Value = 0.00340486317665188

# llogis cov
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0105701973074859
$ws.Range("B2").Value = 0.00506328286632033
$ws.Range("A3").Value = 0.00506328286632033
$ws.Range("B3").Value = 0.0165419987109093

# gompertz cov
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0175035827891505
$ws.Range("B2").Value = -0.00054132285507576
$ws.Range("A3").Value = -0.00054132285507576
$ws.Range("B3").Value = 0.0000454237753083905
